$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.790.99'
$ws.Range("E2").Value = '  -1.69%  '
$ws.Range("D3").Value = '2.225.88'
$ws.Range("E3").Value = '  -1.17%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '''250.93'
$ws.Range("E5").Value = '  +7.26%  '
$ws.Range("D6").Value = '''0.630'
$ws.Range("E6").Value = '  -0.71%  '
$ws.Range("D7").Value = '''72.00'
$ws.Range("E7").Value = '  +3.01%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = '''0.593'
$ws.Range("E9").Value = '  +5.49%  '
$ws.Range("D10").Value = '''41.27'
$ws.Range("E10").Value = '  +13.84%  '
$ws.Range("D11").Value = '''0.0969'
$ws.Range("E11").Value = '  -2.77%  '
$ws.Range("D12").Value = '''58.21'
$ws.Range("E12").Value = '  -0.68%  '
$ws.Range("E13").Value = '  +5.96%  '
$ws.Range("D14").Value = '''0.106'
$ws.Range("E14").Value = '  -0.56%  '
$ws.Range("D15").Value = '2.557.53'
$ws.Range("E15").Value = '  -1.16%  '
$ws.Range("D16").Value = '''15.00'
$ws.Range("E16").Value = '  -1.03%  '
$ws.Range("D17").Value = '''0.867'
$ws.Range("E17").Value = '  +0.67%  '
$ws.Range("D18").Value = '2.230.34'
$ws.Range("E18").Value = '  -0.63%  '
$ws.Range("D19").Value = '41.755.86'
$ws.Range("E19").Value = '  -1.39%  '
$ws.Range("D20").Value = '0.0₃0967'
$ws.Range("E20").Value = '  -1.36%  '
$ws.Range("D21").Value = '''6.23'
$ws.Range("E21").Value = '  -0.64%  '
$ws.Range("D22").Value = '''72.94'
$ws.Range("E22").Value = '  -0.84%  '
$ws.Range("D23").Value = '''235.43'
$ws.Range("E23").Value = '  -0.58%  '
$ws.Range("B24").Value = 'WEMIXToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D24").Value = '''4.28'
$ws.Range("E24").Value = '  +16.70%  '
$ws.Range("B25").Value = 'ImmutableX'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D25").Value = '''2.12'
$ws.Range("E25").Value = '  +5.70%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").Value = '''2.54'
$ws.Range("E27").Value = '  +4.86%  '
$ws.Range("D28").Value = '''10.73'
$ws.Range("E28").Value = '  +6.71%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '''2.19'
$ws.Range("E29").Value = '  -0.70%  '
$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").Value = '''171.36'
$ws.Range("E30").Value = '  +1.72%  '
$ws.Range("E31").Value = '  +0.78%  '
$ws.Range("E32").Value = '  +3.08%  '
$ws.Range("E33").Value = '  -1.24%  '
$ws.Range("D34").Value = '''5.59'
$ws.Range("E34").Value = '  +3.20%  '
$ws.Range("D35").Value = '''0.0734'
$ws.Range("E35").Value = '  +0.37%  '
$ws.Range("D36").Value = '''4.73'
$ws.Range("E36").Value = '  +0.26%  '
$ws.Range("D37").Value = '''26.42'
$ws.Range("E37").Value = '  +21.89%  '
$ws.Range("D38").Value = '''4.00'
$ws.Range("E38").Value = '  +9.76%  '
$ws.Range("D39").Value = '''0.0306'
$ws.Range("E39").Value = '  +12.85%  '
$ws.Range("E40").Value = '  +0.47%  '
$ws.Range("D41").Value = '''5.95'
$ws.Range("E41").Value = '  -1.13%  '
$ws.Range("D42").Value = '''66.89'
$ws.Range("E42").Value = '  +1.55%  '
$ws.Range("D43").Value = '''12.08'
$ws.Range("E43").Value = '  +19.40%  '
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").Value = '''0.203'
$ws.Range("E44").Value = '  +5.27%  '
$ws.Range("B45").Value = 'FTXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D45").Value = '''4.92'
$ws.Range("E45").Value = '  -0.97%  '
$ws.Range("D46").Value = '''8.79'
$ws.Range("E46").Value = '  -3.58%  '
$ws.Range("E47").Value = '  -0.38%  '
$ws.Range("D48").Value = '''4.64'
$ws.Range("E48").Value = '  +2.50%  '
$ws.Range("E49").Value = '  -0.07%  '
$ws.Range("E50").Value = '  +6.70%  '
$ws.Range("D51").Value = '''1.20'
$ws.Range("E51").Value = '  +0.92%  '
